$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '35.544.27'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '1.913.32'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '0.704'
$ws.Range("E5").Value = '  +6.24%  '
$ws.Range("D6").Value = '247.39'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '40.82'
$ws.Range("E8").Value = '  -3.51%  '
$ws.Range("E9").Value = '  +2.67%  '
$ws.Range("D10").Value = '52.86'
$ws.Range("E10").Value = '  +7.49%  '
$ws.Range("D11").Value = '0.0737'
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("D12").Value = '0.0992'
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("D13").Value = '2.190.63'
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("D14").Value = '12.72'
$ws.Range("E14").Value = '  +2.64%  '
$ws.Range("D15").Value = '0.718'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").Value = '1.919.73'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = '4.92'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = '35.533.40'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").Value = '73.32'
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  +3.89%  '
$ws.Range("D22").Value = '242.72'
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").Value = '5.09'
$ws.Range("E23").Value = '  +4.92%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '2.32'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = '2.29'
$ws.Range("E26").Value = '  +2.43%  '
$ws.Range("D27").Value = '168.44'
$ws.Range("E27").Value = '  -1.91%  '
$ws.Range("D28").Value = '8.65'
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("D29").Value = '18.88'
$ws.Range("E29").Value = '  +2.88%  '
$ws.Range("E30").Value = '  +2.23%  '
$ws.Range("D31").Value = '4.142.55'
$ws.Range("D32").Value = '4.27'
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("D33").Value = '0.0580'
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("D34").Value = '1.93'
$ws.Range("E34").Value = '  +11.33%  '
$ws.Range("D35").Value = '4.23'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").Value = '0.918'
$ws.Range("E37").Value = '  -6.32%  '
$ws.Range("D38").Value = '1.49'
$ws.Range("E38").Value = '  +11.14%  '
$ws.Range("D39").Value = '2.05'
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("D40").Value = '17.73'
$ws.Range("E40").Value = '  +13.26%  '
$ws.Range("D41").Value = '99.03'
$ws.Range("E41").Value = '  +6.88%  '
$ws.Range("D42").Value = '1.15'
$ws.Range("E42").Value = '  +2.41%  '
$ws.Range("D43").Value = '0.0211'
$ws.Range("E43").Value = '  +2.14%  '
$ws.Range("D44").Value = '0.0647'
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("D45").Value = '1.353.02'
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").Value = '2.47'
$ws.Range("E46").Value = '  +3.00%  '
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").Value = '45.86'
$ws.Range("E49").Value = '  -3.79%  '
$ws.Range("D50").Value = '12.21'
$ws.Range("E50").Value = '  -3.42%  '
$ws.Range("D51").Value = '6.57'
$ws.Range("E51").Value = '  -0.43%  '
